# Update the "partner match" sheet:
#  - The shared "Solvers" value that used to read
#    "None,AIR-INK: Air-Pollution to ink" (row 33, Save the Children)
#    now reads "None,Mycotech" and belongs to row 27
#    (MIT Environmental Solutions Initative (John Fernandez)) instead.
#  - Row 27's Count moves from 0 to 1.
#  - Row 33 reverts to "None" with Count 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 27: MIT Environmental Solutions Initative (John Fernandez)
$ws.Range("B27").Value = "None,Mycotech"
$ws.Range("C27").Value = 1

# Row 33: Save the Children
$ws.Range("B33").Value = "None"
$ws.Range("C33").Value = 0
